$d = $word.ActiveDocument

function Insert-BulletAfter($anchorText, $newText) {
    # Find the anchor paragraph's text and collapse the range to just
    # after it, so a new paragraph can be appended right after it.
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, `
                                $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $anchorText"
    }
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()

    # Move into the freshly created (empty) paragraph that now sits
    # right after the anchor paragraph mark.
    $insertStart = $rng.End + 1
    $rng.SetRange($insertStart, $insertStart)
    $rng.InsertAfter($newText)

    # The new paragraph inherited the previous bullet's list level
    # (ilvl = 1); bring it back to the top bullet level (ilvl = 0),
    # matching the level used for the other numId=2 bullets here.
    $rng.SetRange($insertStart, $insertStart)
    $rng.ListFormat.ListLevelNumber = 1
}

$anchor1 = "Có một button cho xem đáp án. Khi click vào sẽ hiện ra 1 cái panel, hoặc text box chứa nội dung của bài tập viết đó."
$text1 = "Khi bắt đầu vào form chính tả, Long nghĩ không cần cái pic bắt đầu. Như thế không tiện dụng cho lắm. Thay vì vậy, mình load lên màn hình bài học luôn. Như vậy mới giống đặc tả."
Insert-BulletAfter $anchor1 $text1

$anchor2 = $text1
$text2 = "Khi chọn vào Học tập (ở main form), trên cái pnl sẽ hiện ra 5 chữ: Tập đọc, Chính tả, Luyện từ và câu, Kể chuyện, Tập làm văn, được nằm lộn xộn thay vì chữ Tiếng Việt. Long nghĩ như vậy sẽ hay hơn."
Insert-BulletAfter $anchor2 $text2

Write-Output "done"
